$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H58").Value = 971.25
$ws.Range("I58").Value = 295
$ws.Range("K58").Value = 885
$ws.Range("M58").Value = -735
$ws.Range("H70").Value = 25001362
$ws.Range("I70").Value = 28572628
$ws.Range("K70").Value = 85717884
$ws.Range("M70").Value = -85717614
$ws.Range("H73").Value = 25001362
$ws.Range("I73").Value = 28572628
$ws.Range("K73").Value = 85717884
$ws.Range("M73").Value = -85716948
$ws.Range("H74").Value = 5738
$ws.Range("I74").Value = 5738
$ws.Range("K74").Value = 5738
$ws.Range("M74").Value = -4802
$ws.Range("H77").Value = 5738
$ws.Range("I77").Value = 5738
$ws.Range("K77").Value = 28690
$ws.Range("M77").Value = -24010
$ws.Range("H121").Value = 2321.8
$ws.Range("J121").Value = 2321.8
$ws.Range("L121").Value = 6965.400000000001
$ws.Range("N121").Value = -10459.4
$ws.Range("H137").Value = 624579.3
$ws.Range("I137").Value = 852166.9
$ws.Range("J137").Value = 3885.9092
$ws.Range("K137").Value = 2556500.7
$ws.Range("L137").Value = 11657.7276
$ws.Range("M137").Value = -2553950.7
$ws.Range("N137").Value = -16757.7276
$ws.Range("H138").Value = 4075.13
$ws.Range("J138").Value = 4579.926
$ws.Range("L138").Value = 13739.778
$ws.Range("N138").Value = -24019.778

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 4140.278
$ws.Range("I2").Value = 4794.3076
$ws.Range("K2").Value = 4794.3076
$ws.Range("M2").Value = -4681.3076
$ws.Range("H116").Value = 4140.278
$ws.Range("I116").Value = 4794.3076
$ws.Range("K116").Value = 4794.3076
$ws.Range("M116").Value = -2500.3076
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530
$ws.Range("H132").Value = 1729.2985
$ws.Range("I132").Value = 1047.7322
$ws.Range("K132").Value = 3143.1966
$ws.Range("M132").Value = -613.1965999999998

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 4140.278
$ws.Range("I3").Value = 4794.3076
$ws.Range("K3").Value = 4794.3076
$ws.Range("M3").Value = -4680.3076
$ws.Range("H26").Value = 53839.152
$ws.Range("I26").Value = 41591.2
$ws.Range("K26").Value = 41591.2
$ws.Range("M26").Value = -41299.2
$ws.Range("H99").Value = 42788.57
$ws.Range("I99").Value = 84840
$ws.Range("J99").Value = 11250
$ws.Range("K99").Value = 84840
$ws.Range("L99").Value = 11250
$ws.Range("M99").Value = -83342
$ws.Range("N99").Value = -14246
$ws.Range("H134").Value = 2318.0605
$ws.Range("I134").Value = 1235.3334
$ws.Range("K134").Value = 3706.0002
$ws.Range("M134").Value = -1171.0002

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 4550.8
$ws.Range("I31").Value = 2756
$ws.Range("J31").Value = 4999.5
$ws.Range("K31").Value = 2756
$ws.Range("L31").Value = 4999.5
$ws.Range("M31").Value = -2461
$ws.Range("N31").Value = -5589.5
$ws.Range("H34").Value = 4550.8
$ws.Range("I34").Value = 2756
$ws.Range("J34").Value = 4999.5
$ws.Range("K34").Value = 2756
$ws.Range("L34").Value = 4999.5
$ws.Range("M34").Value = -2554
$ws.Range("N34").Value = -5403.5
$ws.Range("H105").Value = 7978.3687
$ws.Range("I105").Value = 10632.833
$ws.Range("K105").Value = 10632.833
$ws.Range("M105").Value = -8885.833000000001
$ws.Range("H109").Value = 59991
$ws.Range("J109").Value = 59991
$ws.Range("L109").Value = 59991
$ws.Range("N109").Value = -62071
$ws.Range("H114").Value = 50554.668
$ws.Range("J114").Value = 50554.668
$ws.Range("L114").Value = 50554.668
$ws.Range("N114").Value = -59232.668
$ws.Range("H132").Value = 6478.222
$ws.Range("I132").Value = 6478.222
$ws.Range("K132").Value = 19434.666
$ws.Range("M132").Value = -16904.666
$ws.Range("H141").Value = 207870
$ws.Range("J141").Value = 220986.38
$ws.Range("L141").Value = 220986.38
$ws.Range("N141").Value = -231346.38

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H80").Value = 100168020
$ws.Range("J80").Value = 210025
$ws.Range("L80").Value = 630075
$ws.Range("N80").Value = -631947
$ws.Range("H83").Value = 100168020
$ws.Range("J83").Value = 210025
$ws.Range("L83").Value = 1890225
$ws.Range("N83").Value = -1899585
$ws.Range("H109").Value = 1436.2858
$ws.Range("I109").Value = 842.3333
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 2526.9999
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -1486.9999
$ws.Range("N109").Value = -17080
$ws.Range("H113").Value = 1128.9286
$ws.Range("J113").Value = 1179.0834
$ws.Range("L113").Value = 3537.2502
$ws.Range("N113").Value = -7877.2502

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H51").Value = 50527.5
$ws.Range("J51").Value = 50527.5
$ws.Range("L51").Value = 50527.5
$ws.Range("N51").Value = -51545.5
$ws.Range("H113").Value = 20816.5
$ws.Range("I113").Value = 29225
$ws.Range("J113").Value = 3999.5
$ws.Range("K113").Value = 29225
$ws.Range("L113").Value = 3999.5
$ws.Range("M113").Value = -27055
$ws.Range("N113").Value = -8339.5
$ws.Range("H126").Value = 16142.5
$ws.Range("I126").Value = 18367
$ws.Range("K126").Value = 55101
$ws.Range("M126").Value = -52631
$ws.Range("H132").Value = 2607.1904
$ws.Range("I132").Value = 2654.6365
$ws.Range("K132").Value = 7963.9095
$ws.Range("M132").Value = -5433.9095

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 2019.9032
$ws.Range("I16").Value = 1965.2
$ws.Range("J16").Value = 2247.8333
$ws.Range("K16").Value = 1965.2
$ws.Range("L16").Value = 2247.8333
$ws.Range("M16").Value = -1795.2
$ws.Range("N16").Value = -2587.8333
$ws.Range("H61").Value = 17418.688
$ws.Range("I61").Value = 1368.7273
$ws.Range("K61").Value = 1368.7273
$ws.Range("M61").Value = -1166.7273
$ws.Range("H92").Value = 58000
$ws.Range("J92").Value = 58000
$ws.Range("L92").Value = 58000
$ws.Range("N92").Value = -62992
$ws.Range("H113").Value = 17418.688
$ws.Range("I113").Value = 1368.7273
$ws.Range("K113").Value = 1368.7273
$ws.Range("M113").Value = 801.2727
$ws.Range("H123").Value = 140000
$ws.Range("J123").Value = 140000
$ws.Range("L123").Value = 140000
$ws.Range("N123").Value = -149800
$ws.Range("H137").Value = 42857.145
$ws.Range("I137").Value = 20000
$ws.Range("K137").Value = 20000
$ws.Range("M137").Value = -14900
$ws.Range("H139").Value = 59000
$ws.Range("J139").Value = 59000
$ws.Range("L139").Value = 59000
$ws.Range("N139").Value = -69280

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H92").Value = 183373330
$ws.Range("J92").Value = 183373330
$ws.Range("L92").Value = 183373330
$ws.Range("N92").Value = -183378322
$ws.Range("H94").Value = 330000000
$ws.Range("J94").Value = 330000000
$ws.Range("L94").Value = 330000000
$ws.Range("N94").Value = -330001802
$ws.Range("H96").Value = 11112831
$ws.Range("I96").Value = 14287640
$ws.Range("K96").Value = 14287640
$ws.Range("M96").Value = -14286267
$ws.Range("H113").Value = 1417.4572
$ws.Range("I113").Value = 891.069
$ws.Range("J113").Value = 3961.6667
$ws.Range("K113").Value = 2673.207
$ws.Range("L113").Value = 11885.0001
$ws.Range("M113").Value = -503.2069999999999
$ws.Range("N113").Value = -16225.0001
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H129").Value = 75000
$ws.Range("J129").Value = 75000
$ws.Range("L129").Value = 75000
$ws.Range("N129").Value = -85000
$ws.Range("H132").Value = 14576.542
$ws.Range("I132").Value = 20642.666
$ws.Range("J132").Value = 4466.3335
$ws.Range("K132").Value = 61927.99800000001
$ws.Range("L132").Value = 13399.0005
$ws.Range("M132").Value = -59397.99800000001
$ws.Range("N132").Value = -18459.0005
